$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data (and hyperlinks) down
$ws.Rows("2:2").Insert()

# Set the new email value (plain text, no hyperlink, no style)
$ws.Range("A2").Value = "avi@testmail.com"

# Update the active selection
$ws.Range("G13").Select()
